{"js": "// Apply the benchmark-stats refresh to the single results table.\n// Each table row holds one logical stat in its lone cell; the last three\n// rows (1-based 44-46 / 0-based 43-45) previously packed ten tab-separated\n// numbers into one run \u2014 the new content replaces that whole run with a\n// single short value.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\n\n// Map of 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"303\",\n  4: \"0.00002\",\n  5: \"0.00024\",\n  7: \"0.00001\",\n  11: \"0.01171\",\n  43: \"100\",\n  44: \"0.01\",\n  45: \"333\",\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = updates[rowIndexStr];\n}\n\nawait context.sync();\n", "ps1": "# Refresh the benchmark-stats results table: a handful of single-number\n# rows get new values, and the three rows that used to carry ten\n# tab-separated numbers in one run (rows 44-46) are collapsed down to a\n# single short value each.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"303\"\n    5  = \"0.00002\"\n    6  = \"0.00024\"\n    8  = \"0.00001\"\n    12 = \"0.01171\"\n    44 = \"100\"\n    45 = \"0.01\"\n    46 = \"333\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
